$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.847479
$ws.Range("N2").Value = 5.542437
$ws.Range("O2").Value = 0.3181373042830636
$ws.Range("P2").Value = 0.3181373042830637
$ws.Range("Q2").Value = 0.8385251469513333
$ws.Range("R2").Value = 7.546726322562
$ws.Range("S2").Value = 0.3181373042830636
$ws.Range("T2").Value = 0.3181373042830637

# Row 3
$ws.Range("N3").Value = 6.300930999999999
$ws.Range("O3").Value = 0.3616750542791174
$ws.Range("P3").Value = 0.3616750542791174
$ws.Range("R3").Value = 8.579511473805999
$ws.Range("S3").Value = 0.3616750542791174
$ws.Range("T3").Value = 0.3616750542791174

# Row 4
$ws.Range("M4").Value = 1.859385666666666
$ws.Range("N4").Value = 5.578156999999999
$ws.Range("O4").Value = 0.320187641437819
$ws.Range("P4").Value = 0.320187641437819
$ws.Range("Q4").Value = 0.8439292892535555
$ws.Range("R4").Value = 7.595363603281998
$ws.Range("S4").Value = 0.320187641437819
$ws.Range("T4").Value = 0.320187641437819
